$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content change: "Assert" step renamed to "Then" (bdd syntax) ---
$ws.Range("A7").Value = "Then"

# --- Conditional formatting cleanup ---
# The sheet had two duplicate rule-sets: one scoped to "A1:XFD4 A14:XFD1048576 E5:XFD13"
# and a near-identical one scoped to "A5:D13". Collapse down to a single rule-set
# (the A5:D13-derived one) that now applies to the whole sheet, matching the
# consolidated conditional formatting Excel produced after the edit.
$fcs = $ws.Cells.FormatConditions

# Delete the first 13 rules (the "A1:XFD4 A14:XFD1048576 E5:XFD13" duplicate set).
for ($i = 1; $i -le 13; $i++) {
    $fcs.Item(1).Delete()
}

# The remaining 13 rules (formerly "A5:D13") are grouped into 4 conditionalFormatting
# blocks (9 + 1 + 2 + 1 rules). Re-point each block to the whole sheet.
$fullRange = $ws.Range("A1:XFD1048576")
$fcs.Item(1).ModifyAppliesToRange($fullRange)
$fcs.Item(10).ModifyAppliesToRange($fullRange)
$fcs.Item(11).ModifyAppliesToRange($fullRange)
$fcs.Item(13).ModifyAppliesToRange($fullRange)

# Update rule text/formulas: "Assert" -> "Then", and the relative formulas that
# were anchored on A5/A4/A6/XFD5 now anchor on A1/A1048576/A2/XFD1 (new top-left).
$fcs.Item(1).Formula1 = '="Then"'
$fcs.Item(5).Formula1 = '=LEFT(A1,LEN("With Properties"))="With Properties"'
$fcs.Item(6).Formula1 = '=RIGHT(A1,LEN(" table of"))=" table of"'
$fcs.Item(7).Formula1 = '=RIGHT(A1,LEN(" of"))=" of"'
$fcs.Item(8).Formula1 = '=AND((RIGHT(A1048576, 3) = " of"), A2 = "With Properties")'
$fcs.Item(9).Formula1 = '=AND(RIGHT(XFD1, 3) = " of", A2 = "With Properties")'
$fcs.Item(10).Formula1 = '=LEN(TRIM(A1))>0'

# --- Clear the stray cell selection left over from editing (D13) ---
$ws.Range("A1").Select()
